# Update computed power-flow line results for the 380 kV case.
# Values for columns B, C, E, F, G, K, L, M, N across data rows 2-25
# are refreshed with the latest simulation output (columns A, D, H, I, J, O unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.740247270146568
$ws.Range("C2").Value = 0.03229379708800195
$ws.Range("E2").Value = 0.05337578754117889
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.002630005008481243
$ws.Range("K2").Value = 1.258426602509331
$ws.Range("L2").Value = 0.1986043025361823
$ws.Range("M2").Value = 0.3529509103836119
$ws.Range("N2").Value = 4.377671302972288
$ws.Range("B3").Value = 1.703661558086367
$ws.Range("C3").Value = 0.02810134703410938
$ws.Range("E3").Value = 0.0534509493310722
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.00263437294287609
$ws.Range("K3").Value = 1.220431959008096
$ws.Range("L3").Value = 0.1963722079547878
$ws.Range("M3").Value = 0.3465977830136495
$ws.Range("N3").Value = 4.378939327854553
$ws.Range("B4").Value = 1.682149267224673
$ws.Range("C4").Value = 0.02552521028550814
$ws.Range("E4").Value = 0.0535119162904385
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.002637196045817435
$ws.Range("K4").Value = 1.197836834602839
$ws.Range("L4").Value = 0.19509846115659
$ws.Range("M4").Value = 0.3428846951674416
$ws.Range("N4").Value = 4.38049367945149
$ws.Range("B5").Value = 1.673622022533863
$ws.Range("C5").Value = 0.02447483648060711
$ws.Range("E5").Value = 0.05354049594508137
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002638382098123832
$ws.Range("K5").Value = 1.188813402046861
$ws.Range("L5").Value = 0.1946037465545487
$ws.Range("M5").Value = 0.3414187938598054
$ws.Range("N5").Value = 4.381321899446519
$ws.Range("B6").Value = 1.672220524663459
$ws.Range("C6").Value = 0.02430038536681423
$ws.Range("E6").Value = 0.05354546743875677
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002638581195530859
$ws.Range("K6").Value = 1.187326190228532
$ws.Range("L6").Value = 0.1945230705472554
$ws.Range("M6").Value = 0.3411782335768834
$ws.Range("N6").Value = 4.381471183002205
$ws.Range("B7").Value = 1.682033297445827
$ws.Range("C7").Value = 0.0255110470257307
$ws.Range("E7").Value = 0.05351228659003304
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002637211897136483
$ws.Range("K7").Value = 1.197714395672762
$ws.Range("L7").Value = 0.1950916906599716
$ws.Range("M7").Value = 0.3428647343503037
$ws.Range("N7").Value = 4.380504060720455
$ws.Range("B8").Value = 1.727434975264742
$ws.Range("C8").Value = 0.03084859345764812
$ws.Range("E8").Value = 0.05339863342415896
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.002631481841572011
$ws.Range("K8").Value = 1.24517359925909
$ws.Range("L8").Value = 0.1978145992954765
$ws.Range("M8").Value = 0.3507213787801362
$ws.Range("N8").Value = 4.377947339852383
$ws.Range("B9").Value = 1.824027659841875
$ws.Range("C9").Value = 0.04130432908559101
$ws.Range("E9").Value = 0.0532929534909794
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002621360088706733
$ws.Range("K9").Value = 1.344080547585037
$ws.Range("L9").Value = 0.2039219482420407
$ws.Range("M9").Value = 0.367619412410285
$ws.Range("N9").Value = 4.379102901108482
$ws.Range("B10").Value = 1.899628402638541
$ws.Range("C10").Value = 0.04898604953862673
$ws.Range("E10").Value = 0.05328629258584172
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002614595884106387
$ws.Range("K10").Value = 1.420341813240839
$ws.Range("L10").Value = 0.2088777613706583
$ws.Range("M10").Value = 0.3809473517334396
$ws.Range("N10").Value = 4.383734573809562
$ws.Range("B11").Value = 1.935033640230188
$ws.Range("C11").Value = 0.05248213210860797
$ws.Range("E11").Value = 0.05329858169652724
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002611663064774926
$ws.Range("K11").Value = 1.455823652133574
$ws.Range("L11").Value = 0.2112342906941649
$ws.Range("M11").Value = 0.3872098043298706
$ws.Range("N11").Value = 4.386667943163104
$ws.Range("B12").Value = 1.948586814755117
$ws.Range("C12").Value = 0.05380634814351026
$ws.Range("E12").Value = 0.05330542986686737
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002610573105666859
$ws.Range("K12").Value = 1.469373803807883
$ws.Range("L12").Value = 0.2121413326332942
$ws.Range("M12").Value = 0.389609965144345
$ws.Range("N12").Value = 4.387897935578422
$ws.Range("B13").Value = 1.945661401054394
$ws.Range("C13").Value = 0.05352113859446206
$ws.Range("E13").Value = 0.05330385751941868
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002610806931882687
$ws.Range("K13").Value = 1.466450463098568
$ws.Range("L13").Value = 0.2119453323346931
$ws.Range("M13").Value = 0.389091770570765
$ws.Range("N13").Value = 4.387627728036733
$ws.Range("B14").Value = 1.936145741479947
$ws.Range("C14").Value = 0.05259106911219646
$ws.Range("E14").Value = 0.05329910116494752
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002611572980394715
$ws.Range("K14").Value = 1.456936145693419
$ws.Range("L14").Value = 0.2113086195011533
$ws.Range("M14").Value = 0.3874066915898595
$ws.Range("N14").Value = 4.386766744607456
$ws.Range("B15").Value = 1.930336140867212
$ws.Range("C15").Value = 0.05202141979914643
$ws.Range("E15").Value = 0.05329647331017107
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002612044891291007
$ws.Range("K15").Value = 1.451123200748754
$ws.Range("L15").Value = 0.2109205254811428
$ws.Range("M15").Value = 0.3863782709600727
$ws.Range("N15").Value = 4.386254899927991
$ws.Range("B16").Value = 1.897335006080311
$ws.Range("C16").Value = 0.04875761288144531
$ws.Range("E16").Value = 0.05328579701078162
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002614790444229155
$ws.Range("K16").Value = 1.418038919619278
$ws.Range("L16").Value = 0.2087258106230792
$ws.Range("M16").Value = 0.380542101728544
$ws.Range("N16").Value = 4.383559530614377
$ws.Range("B17").Value = 1.877349716794811
$ws.Range("C17").Value = 0.04675585498512191
$ws.Range("E17").Value = 0.05328316510297526
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 0.00261651162166465
$ws.Range("K17").Value = 1.397945399303012
$ws.Range("L17").Value = 0.2074055703597537
$ws.Range("M17").Value = 0.3770129108830957
$ws.Range("N17").Value = 4.382117919503344
$ws.Range("B18").Value = 1.865950147105309
$ws.Range("C18").Value = 0.04560464465848213
$ws.Range("E18").Value = 0.05328309370651674
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002617515181256753
$ws.Range("K18").Value = 1.386462496241194
$ws.Range("L18").Value = 0.2066558140734713
$ws.Range("M18").Value = 0.3750017909752188
$ws.Range("N18").Value = 4.381366505974228
$ws.Range("B19").Value = 1.862106836022178
$ws.Range("C19").Value = 0.04521488725514189
$ws.Range("E19").Value = 0.05328331752294524
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.002617857305832022
$ws.Range("K19").Value = 1.382587344132105
$ws.Range("L19").Value = 0.2064036104050899
$ws.Range("M19").Value = 0.3743240848591114
$ws.Range("N19").Value = 4.381125435422973
$ws.Range("B20").Value = 1.879467306186029
$ws.Range("C20").Value = 0.04696892956614818
$ws.Range("E20").Value = 0.05328329605082693
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002616326994297962
$ws.Range("K20").Value = 1.400076690765502
$ws.Range("L20").Value = 0.2075451175425798
$ws.Range("M20").Value = 0.3773866557755738
$ws.Range("N20").Value = 4.382263330789044
$ws.Range("B21").Value = 1.938936758225623
$ws.Range("C21").Value = 0.05286424350106245
$ws.Range("E21").Value = 0.05330043872804247
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.002611347414236345
$ws.Range("K21").Value = 1.459727636206821
$ws.Range("L21").Value = 0.2114952392607705
$ws.Range("M21").Value = 0.3879008609165027
$ws.Range("N21").Value = 4.387016398704304
$ws.Range("B22").Value = 1.978654442623451
$ws.Range("C22").Value = 0.05671911171830857
$ws.Range("E22").Value = 0.05332443017907096
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002608213202603078
$ws.Range("K22").Value = 1.499377316002381
$ws.Range("L22").Value = 0.2141624090207586
$ws.Range("M22").Value = 0.3949398312316674
$ws.Range("N22").Value = 4.390817696673366
$ws.Range("B23").Value = 1.957378449897078
$ws.Range("C23").Value = 0.05466148819827765
$ws.Range("E23").Value = 0.05331045814125623
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002609875024431632
$ws.Range("K23").Value = 1.478154640906723
$ws.Range("L23").Value = 0.2127310658510595
$ws.Range("M23").Value = 0.3911676834742011
$ws.Range("N23").Value = 4.388725171335977
$ws.Range("B24").Value = 1.878509662950876
$ws.Range("C24").Value = 0.04687259971801438
$ws.Range("E24").Value = 0.05328323235820775
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002616410420734977
$ws.Range("K24").Value = 1.399112918592408
$ws.Range("L24").Value = 0.2074819993945169
$ws.Range("M24").Value = 0.3772176300833081
$ws.Range("N24").Value = 4.382197349301293
$ws.Range("B25").Value = 1.797084843811689
$ws.Range("C25").Value = 0.03847628167778794
$ws.Range("E25").Value = 0.0533090444509714
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002623979704746221
$ws.Range("K25").Value = 1.316694911397008
$ws.Range("L25").Value = 0.2021874865776994
$ws.Range("M25").Value = 0.3628880102616705
$ws.Range("N25").Value = 4.378127600357516
